$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text-like (numeric-looking) value into a cell while keeping it
# stored as literal text (so "0.1420" doesn't turn into 0.142, etc.), then
# restore the cell's style back to Normal so we don't leave stray formatting.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Simple price (column D) corrections - coin/link/other columns unchanged
# ---------------------------------------------------------------------------
$priceUpdates = [ordered]@{
    2  = "245.24"
    4  = "5.404"
    6  = "3.389"
    7  = "0.8073"
    8  = "0.9297"
    9  = "0.1420"
    10 = "0.07437"
    12 = "0.03039"
    13 = "0.09362"
    14 = "3.932"
    15 = "0.001590"
    16 = "0.04835"
    17 = "0.0005943"
    18 = "0.005392"
    19 = "0.004153"
    20 = "0.0009822"
    21 = "0.00007104"
    23 = "6.428"
    40 = "0.03969"
    45 = "0.00005206"
}

foreach ($row in $priceUpdates.Keys) {
    Set-TextValue $ws.Range("D$row") $priceUpdates[$row]
}

# ---------------------------------------------------------------------------
# Rows 41-43 got reshuffled (coins rotated up one slot, with refreshed
# prices), so rewrite B/C/D/E fully for those three rows.
# ---------------------------------------------------------------------------
$rowData = @{
    41 = @{ B = "KickToken";  C = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"; D = "0.006353"; E = "40KickTokenKICK" }
    42 = @{ B = "BKEXToken";  C = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk";       D = "0.1074";   E = "41BKEXTokenBKK" }
    43 = @{ B = "CEJI";       C = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji";           D = "0.002711"; E = "42CEJICEJI" }
}

foreach ($row in $rowData.Keys) {
    $data = $rowData[$row]
    $ws.Range("B$row").Value = $data.B
    $ws.Range("C$row").Value = $data.C
    Set-TextValue $ws.Range("D$row") $data.D
    $ws.Range("E$row").Value = $data.E
}

# ---------------------------------------------------------------------------
# Row 44: price update + "Bestin24h" suffix added to column E
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("D44") "0.006957"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"

# ---------------------------------------------------------------------------
# Row 48: "Bestin24h" suffix removed from column E
# ---------------------------------------------------------------------------
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"
